$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 598
$ws.Range("G4").Value = 536
$ws.Range("G6").Value = 93
$ws.Range("G8").Value = 93
$ws.Range("G9").Value = 17
$ws.Range("G11").Value = 17
$ws.Range("G19").Value = 154
$ws.Range("G23").Value = 78
$ws.Range("G24").Value = 57
$ws.Range("G25").Value = 21

$ws.Range("G2:G25").Select()
